# Inserts a new row of data at row 2 ("南京耀多信息技术有限公司") by shifting the
# existing B:P content of rows 2-11 down into rows 3-12 (the A column, which holds
# a manually-maintained running index, is left untouched row-by-row and only the
# brand-new last row gets an explicit new index value). Finally row 2's B:P cells
# are overwritten with the new company's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the B:P content of existing rows down by one, working from the bottom
#    up so we never overwrite a source row before it has been copied.
for ($r = 11; $r -ge 2; $r--) {
    $src = $ws.Range("B" + $r + ":P" + $r)
    $dst = $ws.Range("B" + ($r + 1) + ":P" + ($r + 1))
    $src.Copy($dst)
}

# 2) The new last row (12) needs its own index in column A (continuing the
#    existing 0-based sequence). Copy the style from the row above, then set
#    the value explicitly.
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A12").Value = 10

# 3) Fill the new row 2 with the new company's details (column A keeps its
#    original value of 0, unchanged).
$ws.Range("B2").Value = "南京耀多信息技术有限公司"
$ws.Range("C2").Value = "江苏南京"
$ws.Range("D2").Value = "技术部"
$ws.Range("E2").Value = "Android"
$ws.Range("F2").Value = "9:00-18:00"
$ws.Range("G2").Value = "1h"
$ws.Range("H2").Value = "一开始996，后来发不起加班费不给加班了，欠的加班费也不发"
$ws.Range("I2").Value = "最低额度"
$ws.Range("J2").Value = "无"
$ws.Range("K2").Value = "八折"
$ws.Range("L2").Value = "提供笔记本"
$ws.Range("M2").Value = "有"
$ws.Range("N2").Value = "钉钉位置打卡"
$ws.Range("O2").Value = "老板阴晴不定，随意开除员工"
$ws.Range("P2").Value = ""
